$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Sina"
$ws.Range("C4").Value = "Rewrite     -   Lily"
